$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), matching the formatting of the existing
# header cells (bold, centered, bordered) by copying G1's format into H1.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data row value for the new column.
$ws.Range("H2").Value = 0
